$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "N/A"
$ws.Range("N2").Value = "N/A"
$ws.Range("O2").Value = "N/A"
